# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# timestamps for the en -> zh-cn / de-de handback rows with the latest
# report-generation run's values.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-17 22:52:43"
$zhcn.Range("H2").Value = "2016-03-17 22:53:07"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-17 22:52:46"
$dede.Range("H2").Value = "2016-03-17 22:53:13"
